# Add three new rows of work-log data to the bottom of the table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: 王康明 - 继续编写前端代码 - 210 lines of code - 15 hours
$ws.Range("A8").Value = "王康明"
$ws.Range("B8").Value = "继续编写前端代码"
$ws.Range("C8").Value = 210
$ws.Range("E8").Value = 15
$ws.Rows.Item(8).RowHeight = 21.5

# Row 9: 王康明 - 根据本小组评审意见修改bug - 4 hours
$ws.Range("A9").Value = "王康明"
$ws.Range("B9").Value = "根据本小组评审意见修改bug"
$ws.Range("E9").Value = 4
$ws.Rows.Item(9).RowHeight = 28

# Row 10: 王康明 - 根据其他小组评审意见修改bug - 0.5 hours
$ws.Range("A10").Value = "王康明"
$ws.Range("B10").Value = "根据其他小组评审意见修改bug"
$ws.Range("E10").Value = 0.5
$ws.Rows.Item(10).RowHeight = 28

# Update the active selection to reflect where editing finished
$ws.Range("H10").Select()
